# Reverse the order of the comma-separated "Recorded By" entries in
# column G (the list of users/systems that recorded each attendance
# session). Cells that contain only a single entry (no comma) are left
# untouched, matching the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G holds "Recorded By"; data starts on row 2 (row 1 is the header).
$startRow = [Math]::Max($firstRow, 2)

for ($r = $startRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Text

    if ($value -ne $null -and $value -like "*,*") {
        $parts = $value -split ", "
        $reversed = $parts[($parts.Length - 1)..0]
        $newValue = [string]::Join(", ", $reversed)
        $cell.Value = $newValue
    }
}
